# "Generate Report for Handoff"
#
# The localization-status report is regenerated right before handoff
# instead of after handback, so:
#   - the per-language status that read "Handed back: in sync with en-US"
#     now reads "Ready for handoff" (Overview!E2/F2 and the "Status"
#     column on each language sheet all share this text);
#   - the report timestamps move forward a couple of minutes to the new
#     generation time;
#   - the Overview / Status columns shrink now that "Ready for handoff"
#     is shorter than the old message, so their widths are refreshed.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Report generation / handoff timestamps.
$wsOverview.Range("G2").Value = "2016-08-30 09:25:11"
$wsDeDe.Range("H2").Value = "2016-08-30 09:25:11"
$wsZhCn.Range("H2").Value = "2016-08-30 09:24:52"

# Column widths refreshed to fit the new, shorter status text.
$wsOverview.Range("E1").ColumnWidth = 17.2159881591797
$wsOverview.Range("F1").ColumnWidth = 17.2159881591797
$wsZhCn.Range("C1").ColumnWidth = 17.2159881591797
$wsDeDe.Range("C1").ColumnWidth = 17.2159881591797
